$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: dimension/measure identifiers
$ws.Range("A3").Value = "sdmx-dimension:refArea"
$ws.Range("B3").Value = "iaest-dimension:entidad-singular"
$ws.Range("D3").Value = "iaest-dimension:nucleo"
$ws.Range("E3").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "sdmx-dimension:refArea"
$ws.Range("G3").Value = "iaest-dimension:nucleodiseminado"

# Row 4: dim/medida classification
$ws.Range("A4").Value = "dim"
$ws.Range("B4").Value = "dim"
$ws.Range("D4").Value = "dim"
$ws.Range("E4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"

# Row 5: type / concept scheme
$ws.Range("A5").Value = "URI-Municipio"
$ws.Range("B5").Value = "skos:Concept"
$ws.Range("D5").Value = "skos:Concept"
$ws.Range("E5").Value = "URI-comarca"
$ws.Range("F5").Value = "URI-Provincia"
$ws.Range("G5").Value = "skos:Concept"

# Row 6 (new): mapping file references
$ws.Range("B6").Value = "mapping-entidad-singular.xlsx"
$ws.Range("D6").Value = "mapping-nucleo.xlsx"
$ws.Range("G6").Value = "mapping-nucleodiseminado.xlsx"

# Copy the row-5 cell style (s="1") onto the new row-6 cells so they match
# the existing sheet formatting.
$ws.Range("A5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
